# edit.ps1
# Commit message: "Updated for TC003(new Scenario)"
#
# The workbook keeps two tabs with the same logical "AssignLeave" test-data
# layout (TestCase, Execute, EmployeeName, LeaveType, ... ValidationMessage).
# This edit:
#   * Promotes "Sheet1" (previously a slimmed-down 9-column summary) to the
#     full 18-column (A:R) layout and adds a brand-new "TC003" scenario row
#     (Vacation for Russel / All Day / Half Day / Afternoon).
#   * Trims the "AssignLeave" tab back down to just the header plus the
#     TC001 and TC002 rows, and corrects TC001s employee name to "Hannah".

$wb = $excel.ActiveWorkbook
$wsAssign = $wb.Worksheets.Item("AssignLeave")
$wsData = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# "AssignLeave": clear everything, then rewrite header + TC001 + TC002
# ------------------------------------------------------------------
$wsAssign.Cells.Clear()

# Row 1: header
$wsAssign.Range("A1").Value = "Execute"
$wsAssign.Range("B1").Value = "TestCase"
$wsAssign.Range("C1").Value = "EmployeeName"
$wsAssign.Range("D1").Value = "LeaveType"
$wsAssign.Range("E1").Value = "FromDate"
$wsAssign.Range("F1").Value = "ToDate"
$wsAssign.Range("G1").Value = "PartDays"
$wsAssign.Range("H1").Value = "FirstDuration"
$wsAssign.Range("I1").Value = "FirstAMPM"
$wsAssign.Range("J1").Value = "StartDayFrom"
$wsAssign.Range("K1").Value = "StartDatTo"
$wsAssign.Range("L1").Value = "SecondDuration"
$wsAssign.Range("M1").Value = "SecondAMPM"
$wsAssign.Range("N1").Value = "EndDayFrom"
$wsAssign.Range("O1").Value = "EndDayTo"
$wsAssign.Range("P1").Value = "Comments"
$wsAssign.Range("Q1").Value = "LeaveBalance"
$wsAssign.Range("R1").Value = "ValidationMessage"

# Row 2: TC001
$wsAssign.Range("A2").Value = "Y"
$wsAssign.Range("B2").Value = "TC001"
$wsAssign.Range("C2").Value = "Hannah"
$wsAssign.Range("D2").Value = "Maternity US"
$wsAssign.Range("E2").Value = "'2017-09-15"
$wsAssign.Range("F2").Value = "'2017-10-15"
$wsAssign.Range("G2").Value = "None"
$wsAssign.Range("P2").Value = "Maternity leave for Thomos"
$wsAssign.Range("Q2").Value = "'-20"

# Row 3: TC002
$wsAssign.Range("A3").Value = "Y"
$wsAssign.Range("B3").Value = "TC002"
$wsAssign.Range("C3").Value = "John"
$wsAssign.Range("D3").Value = "Vacation US"
$wsAssign.Range("E3").Value = "'2017-09-18"
$wsAssign.Range("F3").Value = "'2017-09-17"
$wsAssign.Range("G3").Value = "None"
$wsAssign.Range("R3").Value = "To date should be after from date"

$wsAssign.Activate()
$wsAssign.Range("K13").Select()

# ------------------------------------------------------------------
# "Sheet1": clear everything, then rewrite full layout incl. new TC003
# ------------------------------------------------------------------
$wsData.Cells.Clear()

# Row 1: header
$wsData.Range("A1").Value = "Execute"
$wsData.Range("B1").Value = "TestCase"
$wsData.Range("C1").Value = "EmployeeName"
$wsData.Range("D1").Value = "LeaveType"
$wsData.Range("E1").Value = "FromDate"
$wsData.Range("F1").Value = "ToDate"
$wsData.Range("G1").Value = "PartDays"
$wsData.Range("H1").Value = "FirstDuration"
$wsData.Range("I1").Value = "FirstAMPM"
$wsData.Range("J1").Value = "StartDayFrom"
$wsData.Range("K1").Value = "StartDatTo"
$wsData.Range("L1").Value = "SecondDuration"
$wsData.Range("M1").Value = "SecondAMPM"
$wsData.Range("N1").Value = "EndDayFrom"
$wsData.Range("O1").Value = "EndDayTo"
$wsData.Range("P1").Value = "Comments"
$wsData.Range("Q1").Value = "LeaveBalance"
$wsData.Range("R1").Value = "ValidationMessage"

# Row 2: TC001
$wsData.Range("A2").Value = "Y"
$wsData.Range("B2").Value = "TC001"
$wsData.Range("C2").Value = "Hannah"
$wsData.Range("D2").Value = "Maternity US"
$wsData.Range("E2").Value = "'2017-09-15"
$wsData.Range("F2").Value = "'2017-10-15"
$wsData.Range("G2").Value = "None"
$wsData.Range("P2").Value = "Maternity leave for Thomos"
$wsData.Range("Q2").Value = "'-20"

# Row 3: TC002
$wsData.Range("A3").Value = "Y"
$wsData.Range("B3").Value = "TC002"
$wsData.Range("C3").Value = "John"
$wsData.Range("D3").Value = "Vacation US"
$wsData.Range("E3").Value = "'2017-09-18"
$wsData.Range("F3").Value = "'2017-09-17"
$wsData.Range("G3").Value = "None"
$wsData.Range("R3").Value = "To date should be after from date"

# Row 4: TC003 (new scenario)
$wsData.Range("A4").Value = "Y"
$wsData.Range("B4").Value = "TC003"
$wsData.Range("C4").Value = "John"
$wsData.Range("D4").Value = "Vacation US"
$wsData.Range("E4").Value = "'2017-09-15"
$wsData.Range("F4").Value = "'2017-09-17"
$wsData.Range("G4").Value = "All Day"
$wsData.Range("H4").Value = "Half Day"
$wsData.Range("I4").Value = "Afternoon"
$wsData.Range("P4").Value = "Vacation for Russel"
$wsData.Range("Q4").Value = "'-2"

# Row 5: TC004
$wsData.Range("N5").NumberFormat = "h:mm"
$wsData.Range("A5").Value = "Y"
$wsData.Range("B5").Value = "TC004"
$wsData.Range("C5").Value = "Fiona"
$wsData.Range("D5").Value = "FMLA US"
$wsData.Range("E5").Value = "'2017-09-15"
$wsData.Range("F5").Value = "'2017-09-17"
$wsData.Range("G5").Value = "Start and End Day"
$wsData.Range("H5").Value = "Half Day"
$wsData.Range("I5").Value = "Afternoon"
$wsData.Range("L5").Value = "Specify Time"
$wsData.Range("N5").Value = "'04:00"
$wsData.Range("O5").Value = "'17:00"
$wsData.Range("R5").Value = "Duration should be less than work shift length"

$wsData.Activate()
$wsData.Range("A1:XFD5").Select()

$wsAssign.Activate()
